$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.844.47'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '1.840.55'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '231.63'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.619'
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '40.01'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.329'
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0982'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').Value = '2.107.65'
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.51'
$ws.Range('E13').Value = '  +4.31%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.842.66'
$ws.Range('E14').Value = '  +1.09%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.673'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').Value = '34.857.68'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.85'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').Value = '0.0₃0789'
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '240.95'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('E21').Value = '  +2.40%  '
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '171.31'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.80'
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.47'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  +2.30%  '
$ws.Range('E29').Value = '  -3.75%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('E32').Value = '  -4.63%  '
$ws.Range('E33').Value = '  -1.74%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.91'
$ws.Range('E34').Value = '  +7.73%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.23'
$ws.Range('E35').Value = '  +7.60%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.45'
$ws.Range('E36').Value = '  +13.08%  '
$ws.Range('E37').Value = '  +1.23%  '
$ws.Range('E38').Value = '  +7.05%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '90.24'
$ws.Range('E39').Value = '  -1.70%  '
$ws.Range('D40').Value = '1.344.89'
$ws.Range('E40').Value = '  +2.54%  '
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '14.91'
$ws.Range('E42').Value = '  +3.72%  '
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('E44').Value = '  -2.88%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.75'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('E47').Value = '  +2.07%  '
$ws.Range('D48').Value = '2.021.00'
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('E49').Value = '  +22.47%  '
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0666'
$ws.Range('E51').Value = '  +1.18%  '
